# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "G" column (header "K") held values that were computed from an older
# "Strike#" based formula. This re-derives/regenerates those values with the
# corrected ("K") calculation and writes the new numbers back into column G
# for every data row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new K (column G) value
$kValues = @{
    2  = 4
    3  = 0
    4  = 1
    5  = 2
    6  = 5
    7  = 5
    8  = 1
    9  = 0
    10 = 3
    11 = 2
    12 = 2
    13 = 0
    14 = 1
    15 = 2
    17 = 0
    18 = 3
    19 = 0
    20 = 0
    21 = 2
    22 = 0
    23 = 2
    24 = 0
    25 = 2
    26 = 3
    27 = 1
    28 = 0
    29 = 1
    30 = 0
    31 = 1
    32 = 1
    33 = 1
    34 = 2
    35 = 1
    36 = 2
    37 = 0
    38 = 2
    39 = 2
    40 = 2
    42 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
